$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(76,0,24,0,0,0,52),
    @(76,0,31,12,4,4,25),
    @(76,0,36,13,4,4,19),
    @(76,0,35,5,7,7,22),
    @(76,0,30,10,2,11,23),
    @(76,0,31,10,4,7,24),
    @(76,0,33,12,5,3,23),
    @(76,0,30,13,4,4,25),
    @(76,0,29,6,6,3,32),
    @(76,0,33,16,1,7,19),
    @(76,0,29,14,5,7,21),
    @(76,0,75,1,0,0,0),
    @(76,0,35,11,4,8,18),
    @(76,0,38,4,5,3,26),
    @(76,0,31,7,4,4,30),
    @(76,0,28,6,6,7,29),
    @(76,0,30,10,3,9,24),
    @(76,0,76,0,0,0,0),
    @(76,0,36,6,1,4,29),
    @(76,0,30,10,8,5,23),
    @(76,0,33,15,1,7,20),
    @(76,0,31,12,3,6,24),
    @(76,0,34,19,4,6,13),
    @(76,0,35,9,7,6,19),
    @(76,0,30,8,2,7,29),
    @(76,0,27,7,5,5,32),
    @(76,0,29,10,2,6,29),
    @(76,0,31,6,2,8,29),
    @(76,0,29,13,5,7,22),
    @(76,0,35,7,8,6,20),
    @(76,0,31,12,2,9,22)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}
